$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the header cells: insert WIN, TOP2, TOP4, RELEGATION before the existing ExpPoints column
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"
$ws.Range("C1:G1").Style = $ws.Range("A1").Style

for ($row = 2; $row -le 19; $row++) {
    $expPoints = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = $expPoints
}
